$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 2) with the values "5000.0" in columns A and B,
# stored as text (matching the original inline-string cells) and with
# the default (unstyled) formatting used by the rest of the sheet.
$ws.Range("A2").Value = "'5000.0"
$ws.Range("B2").Value = "'5000.0"

$ws.Range("A2:B2").Style = "Normal"
